$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Add Sheet2 after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)

# New shared strings must land in this order to match the target file
# (Andrej Karpathy, Langchain, 5 day AI Agents, Illustrated DSA Book,
#  75 blind, then the Sheet1 backprop note last).
$ws2.Range("C5").Value = "Andrej Karpathy"
$ws2.Range("C6").Value = "Langchain"
$ws2.Range("C7").Value = "5 day AI Agents"
$ws2.Range("C8").Value = "Illustrated DSA Book"
$ws2.Range("C9").Value = "75 blind"

$ws2.Columns.Item(3).ColumnWidth = 19

# --- Sheet1: add row 9 (Automatic Backpropogation through _backward(), tanh) ---
$ws1.Range("A9").Value = 3
$ws1.Range("B9").Value = "Automatic Backpropogation through _backward, tanh "

# --- Selections to match final saved state ---
[void]$ws2.Range("H11").Select()

[void]$ws1.Select()
[void]$ws1.Range("B12").Select()
